# Apply updated distance/size labels across the order sheet.
# Rename mapping:
#   D51 -> D55
#   D64 -> D69
#   D80 -> D86
#   S30 -> S31
# (S20 / S25 stay the same). This affects the Condition, Filename_Left,
# Filename_Right, Distance and Size columns (every textual cell that embeds
# one of the Dxx / Sxx tokens).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$lastCol = $usedRange.Columns.Count

# Find header row (row 1) column indexes for the columns we need to touch.
$targetHeaders = @("Condition", "Filename_Left", "Filename_Right", "Distance", "Size")
$colIndexes = @()
for ($c = 1; $c -le $lastCol; $c++) {
    $header = $ws.Cells.Item(1, $c).Value()
    if ($targetHeaders -contains $header) {
        $colIndexes += $c
    }
}

for ($i = 0; $i -lt $colIndexes.Count; $i++) {
    $c = $colIndexes[$i]
    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value()
        if ($val -ne $null -and $val -is [string]) {
            $newVal = $val.Replace("D51", "D55").Replace("D64", "D69").Replace("D80", "D86").Replace("S30", "S31")
            if ($newVal -ne $val) {
                $cell.Value = $newVal
            }
        }
    }
}
